$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "Version 3"
$ws.Range("A1:A3").Select() | Out-Null
